# Update cryptocurrency price/volume data on the worksheet.
# This mirrors an automated "Updated cryptos list" GitHub Actions run,
# refreshing the Price (column D) and Volume(1h) (column E) values for
# each coin row, and correcting the order/data of the Stacks / Stellar
# rows (41 and 42) which had been swapped.
#
# NOTE: several Price values look like plain numbers (e.g. "378.70",
# "0.0835"). Assigning those strings straight to Range.Value would let
# Excel auto-convert them to real numbers (dropping significant trailing
# zeros / introducing floating point noise), whereas the source keeps
# them as literal text. For those cells we set the value with a leading
# apostrophe (Excel's standard "force text" convention) and then reset
# the cell style back to Normal so no stray number-format/quote-prefix
# style lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.035.63"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.947.79"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Formula = "'378.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").Formula = "'101.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("E7").Value = "  -1.14%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Formula = "'36.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.10%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Formula = "'0.0835"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "3.408.10"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").Formula = "'17.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "2.931.00"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Formula = "'0.980"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.53%  "
$ws.Range("D18").Value = "50.952.30"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("E19").Value = "  -6.88%  "
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").Formula = "'12.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.68%  "
$ws.Range("D22").Value = "0.0₃0948"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").Formula = "'68.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Formula = "'260.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("E25").Value = "  +3.59%  "
$ws.Range("D26").Formula = "'8.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.34%  "
$ws.Range("D27").Formula = "'7.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.44%  "
$ws.Range("D28").Formula = "'4.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +10.57%  "
$ws.Range("E31").Value = "  -2.64%  "
$ws.Range("D32").Formula = "'25.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").Formula = "'9.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").Formula = "'50.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.28%  "
$ws.Range("E35").Value = "  -2.98%  "
$ws.Range("D36").Formula = "'33.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("D37").Formula = "'0.0441"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.51%  "
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("E39").Value = "  -1.71%  "
$ws.Range("D40").Formula = "'16.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Formula = "'0.114"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Formula = "'2.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("E43").Value = "  -2.89%  "
$ws.Range("D44").Formula = "'121.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Formula = "'20.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.96%  "
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").Formula = "'0.272"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D49").Value = "1.997.01"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("D51").Formula = "'0.0329"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.63%  "
